# Applies the TENNESSEE_2015.xlsx cleanup edit:
#  1. Rename header row (row 1) from Spanish labels to short code-style column names.
#  2. Clean up the stray "_x000D_"/newline artifact in A222 and title-case it.
#  3. Title-case the Spanish connector words (de/del/el/la/los/las/y) that appear
#     lower-cased inside the state (col A) / municipality (col B) names.
#  4. Re-write a handful of percentage values in column D whose last-bit floating
#     point representation changed.
#  5. Delete the trailing metadata/footnote rows (1483-1487).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# 1. Header row
# ---------------------------------------------------------------------------
$ws.Range("A1").Value2 = "mx_state"
$ws.Range("B1").Value2 = "mx_municipality"
$ws.Range("C1").Value2 = "n_matriculas"
$ws.Range("D1").Value2 = "pct_matriculas"

# ---------------------------------------------------------------------------
# 2. A222 special case: drop the stray "_x000D_"/newline artifact and fix casing
# ---------------------------------------------------------------------------
$ws.Range("A222").Value2 = "Estado De México"

# ---------------------------------------------------------------------------
# 3. Title-case connector words in columns A and B, rows 2..1481
# ---------------------------------------------------------------------------
function Fix-Connectors($text) {
    $result = $text
    $result = $result -replace '\bde\b', 'De'
    $result = $result -replace '\bdel\b', 'Del'
    $result = $result -replace '\bel\b', 'El'
    $result = $result -replace '\bla\b', 'La'
    $result = $result -replace '\blos\b', 'Los'
    $result = $result -replace '\blas\b', 'Las'
    $result = $result -replace '\by\b', 'Y'
    return $result
}

### NOTE: this runtime's `-eq`/`-ne` string comparisons are case-insensitive
### (even with the `-c` prefix), so we cannot reliably detect "did the casing
### change" that way. Instead we always write back the regex-transformed
### text; re-applying the substitutions to text that needs no change is a
### harmless no-op.
$lastRow = 1481
for ($r = 2; $r -le $lastRow; $r++) {
    foreach ($col in @("A", "B")) {
        $cell = $ws.Range("$col$r")
        $val = $cell.Value2
        if ($val -ne $null -and $val.GetType().Name -eq "String") {
            $cell.Value2 = Fix-Connectors $val
        }
    }
}

# ---------------------------------------------------------------------------
# 4. Floating point last-bit corrections in column D
# ---------------------------------------------------------------------------
$rowsToFix1 = @(198, 359, 363, 419, 424, 549, 728, 1126, 1235, 1367, 1390, 1464)
foreach ($r in $rowsToFix1) {
    $ws.Range("D$r").Value2 = 0.0009887529353602767
}
$ws.Range("D369").Value2 = 0.009702138178222715

# ---------------------------------------------------------------------------
# 5. Delete the trailing metadata rows (1483-1487)
# ---------------------------------------------------------------------------
$ws.Range("A1483:A1487").EntireRow.Delete()
